# Clean up non-used exams from the questions/points sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row for chisquare-06.Rnw (row 41 in the original sheet).
$ws.Rows.Item(41).Delete()

# Update points values for several linreg exams (rows shift up by 1 after the
# previous delete: originally 45,47,48,49,50 -> now 44,46,47,48,49).
$ws.Cells.Item(44, 2).Value = 9   # linreg-04.Rnw: 8 -> 9
$ws.Cells.Item(46, 2).Value = 12  # linreg-06.Rnw: 9 -> 12
$ws.Cells.Item(47, 2).Value = 10  # linreg-07.Rnw: 12 -> 10
$ws.Cells.Item(48, 2).Value = 12  # linreg-08.Rnw: 6 -> 12
$ws.Cells.Item(49, 2).Value = 9   # linreg-09.Rnw: 7 -> 9

# Delete the rows for linreg-10.Rnw, linreg-11.Rnw, linreg-12.Rnw
# (originally 51,52,53 -> now 50,51,52 after the first delete).
$ws.Rows.Item(50).Delete()
$ws.Rows.Item(50).Delete()
$ws.Rows.Item(50).Delete()

# Shrink the conditional formatting range to match the new data extent.
$fc = $ws.Range("B2:B83").FormatConditions.Item(1)
$fc.ModifyAppliesToRange($ws.Range("B2:B79"))

# Restore the active selection to match where the author ended up editing.
$ws.Range("A49").Select()
